$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute(
    "havainnointijaksot vuonna Leijonan tähtikuvio 2022: 14.-23.4., 14.-23.5",
    $false,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "Leijonan tähtikuvio havainnointijaksot vuonna 2022: 14.-23.4., 14.-23.5",
    2
)
